$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) retains text formatting so numeric-looking
# strings (e.g. "6.200", "28.001.82") are not coerced into numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.001.82"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "1.861.79"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "312.32"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "0.5123"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").Value = "0.3846"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").Value = "0.08254"
$ws.Range("E9").Value = "  -8.53%  "
$ws.Range("D10").Value = "1.109"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").Value = "41.52"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "6.200"
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("D13").Value = "20.58"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").Value = "1.868.02"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "0.00001097"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").Value = "90.65"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "0.06653"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").Value = "17.67"
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "6.008"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").Value = "28.032.14"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("D25").Value = "2.256"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "2.077.07"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").Value = "2.518"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "157.95"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").Value = "20.46"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "124.65"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("D31").Value = "0.1064"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("D33").Value = "5.952"
$ws.Range("E33").Value = "  +5.97%  "
$ws.Range("D34").Value = "3.594"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "9.331"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").Value = "0.02415"
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").Value = "0.06485"
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").Value = "0.2169"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").Value = "0.6562"
$ws.Range("E39").Value = "  +2.38%  "
$ws.Range("D40").Value = "1.195"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").Value = "4.989"
$ws.Range("E41").Value = "  +1.31%  "
$ws.Range("D42").Value = "1.223"
$ws.Range("E42").Value = "  -4.80%  "
$ws.Range("D43").Value = "11.16"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "0.6158"
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("D45").Value = "12.99"
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("D46").Value = "1.280"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "3.657"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "2.009"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").Value = "1.214"
$ws.Range("E49").Value = "  -2.46%  "
$ws.Range("D50").Value = "120.21"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").Value = "78.61"
$ws.Range("E51").Value = "  -1.13%  "

# Restore the default (no explicit number format) style so the cell
# styling matches the original workbook while keeping the values as text.
$ws.Range("D2:D51").Style = "Normal"
